$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. Numeric-looking text (Price/Volume columns D & E) is
# written with a leading apostrophe so Excel keeps it as text instead of
# coercing it to a number/date, then the style is reset to "Normal" so no
# stray number-format is left on the cell.
$updates = @(
    @{ Cell = "D2"; Value = "26.879.08" }
    @{ Cell = "E2"; Value = "  -0.20%  " }
    @{ Cell = "D3"; Value = "1.861.60" }
    @{ Cell = "E3"; Value = "  +0.11%  " }
    @{ Cell = "D4"; Value = "0.9999" }
    @{ Cell = "E4"; Value = "  -0.16%  " }
    @{ Cell = "D5"; Value = "304.76" }
    @{ Cell = "E5"; Value = "  -0.30%  " }
    @{ Cell = "D6"; Value = "0.9998" }
    @{ Cell = "E6"; Value = "  -0.08%  " }
    @{ Cell = "D7"; Value = "0.5041" }
    @{ Cell = "E7"; Value = "  -0.18%  " }
    @{ Cell = "D8"; Value = "0.3645" }
    @{ Cell = "E8"; Value = "  -2.36%  " }
    @{ Cell = "D9"; Value = "0.07166" }
    @{ Cell = "E9"; Value = "  +0.47%  " }
    @{ Cell = "D10"; Value = "0.8918" }
    @{ Cell = "E10"; Value = "  +0.60%  " }
    @{ Cell = "D11"; Value = "20.67" }
    @{ Cell = "E11"; Value = "  +0.88%  " }
    @{ Cell = "B12"; Value = "TRON" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx" }
    @{ Cell = "D12"; Value = "0.07501" }
    @{ Cell = "E12"; Value = "  -0.63%  " }
    @{ Cell = "B13"; Value = "WrappedEther" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" }
    @{ Cell = "D13"; Value = "1.858.73" }
    @{ Cell = "E13"; Value = "  -0.03%  " }
    @{ Cell = "D14"; Value = "94.94" }
    @{ Cell = "E14"; Value = "  +6.75%  " }
    @{ Cell = "D15"; Value = "5.228" }
    @{ Cell = "E15"; Value = "  -1.09%  " }
    @{ Cell = "D16"; Value = "1.000" }
    @{ Cell = "E16"; Value = "  -0.20%  " }
    @{ Cell = "D17"; Value = "0.000008510" }
    @{ Cell = "E17"; Value = "  +1.86%  " }
    @{ Cell = "D18"; Value = "14.20" }
    @{ Cell = "E18"; Value = "  +1.19%  " }
    @{ Cell = "D20"; Value = "26.936.30" }
    @{ Cell = "E20"; Value = "  -0.17%  " }
    @{ Cell = "D21"; Value = "5.026" }
    @{ Cell = "E21"; Value = "  -0.34%  " }
    @{ Cell = "D22"; Value = "2.105.83" }
    @{ Cell = "E22"; Value = "  +0.73%  " }
    @{ Cell = "D23"; Value = "10.38" }
    @{ Cell = "E23"; Value = "  -0.92%  " }
    @{ Cell = "D24"; Value = "6.401" }
    @{ Cell = "E24"; Value = "  -0.87%  " }
    @{ Cell = "D25"; Value = "147.80" }
    @{ Cell = "E25"; Value = "  +0.37%  " }
    @{ Cell = "D26"; Value = "1.783" }
    @{ Cell = "E26"; Value = "  -3.67%  " }
    @{ Cell = "D27"; Value = "17.87" }
    @{ Cell = "E27"; Value = "  -0.42%  " }
    @{ Cell = "D28"; Value = "2.084" }
    @{ Cell = "E28"; Value = "  -0.29%  " }
    @{ Cell = "D29"; Value = "113.26" }
    @{ Cell = "E29"; Value = "  +0.59%  " }
    @{ Cell = "D30"; Value = "4.707" }
    @{ Cell = "D31"; Value = "4.666" }
    @{ Cell = "E31"; Value = "  +0.71%  " }
    @{ Cell = "D32"; Value = "0.09219" }
    @{ Cell = "E32"; Value = "  +2.07%  " }
    @{ Cell = "D33"; Value = "0.05142" }
    @{ Cell = "E33"; Value = "  +0.69%  " }
    @{ Cell = "D34"; Value = "0.7499" }
    @{ Cell = "E34"; Value = "  +3.16%  " }
    @{ Cell = "D35"; Value = "2.950" }
    @{ Cell = "E35"; Value = "  -3.32%  " }
    @{ Cell = "D36"; Value = "1.153" }
    @{ Cell = "E36"; Value = "  +0.46%  " }
    @{ Cell = "D37"; Value = "3.258" }
    @{ Cell = "E37"; Value = "  +7.30%  " }
    @{ Cell = "D38"; Value = "2.585" }
    @{ Cell = "E38"; Value = "  +5.54%  " }
    @{ Cell = "D39"; Value = "0.02000" }
    @{ Cell = "E39"; Value = "  -1.73%  " }
    @{ Cell = "D40"; Value = "0.5567" }
    @{ Cell = "E40"; Value = "  +4.91%  " }
    @{ Cell = "D41"; Value = "1.069" }
    @{ Cell = "E41"; Value = "  -0.01%  " }
    @{ Cell = "D42"; Value = "6.554" }
    @{ Cell = "E42"; Value = "  -0.18%  " }
    @{ Cell = "D43"; Value = "116.41" }
    @{ Cell = "E43"; Value = "  +1.08%  " }
    @{ Cell = "D44"; Value = "8.585" }
    @{ Cell = "E44"; Value = "  +3.83%  " }
    @{ Cell = "D45"; Value = "0.1471" }
    @{ Cell = "E45"; Value = "  +0.14%  " }
    @{ Cell = "D46"; Value = "0.4694" }
    @{ Cell = "E46"; Value = "  +1.98%  " }
    @{ Cell = "D47"; Value = "0.9992" }
    @{ Cell = "E47"; Value = "  -0.11%  " }
    @{ Cell = "D48"; Value = "10.02" }
    @{ Cell = "E48"; Value = "  +0.24%  " }
    @{ Cell = "E49"; Value = "  +0.44%  " }
    @{ Cell = "D50"; Value = "36.69" }
    @{ Cell = "E50"; Value = "  +0.60%  " }
    @{ Cell = "D51"; Value = "62.97" }
    @{ Cell = "E51"; Value = "  -1.41%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $col = $u.Cell.Substring(0, 1)
    if ($col -eq "D" -or $col -eq "E") {
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
